# The author added a new weekly record for "Jengibre" at Vega Modelo de
# Temuco, inserted as row 147 (the existing rows 147-256 shift down to
# 148-257, and the sheet's used range grows from A1:R256 to A1:R257).
#
# The new row carries the same market/category/quality/unit/origin data as
# the (old) row 147 it was inserted in front of, but with its own date and
# volume:
#   D (Fecha)   -> 44978  (old row147 had 44972)
#   J (Volumen) -> 15     (old row147 had 30)
# Columns K/L/M/P etc. keep the same values the old row147 had.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 147..256 down to 148..257, leaving a blank row 147 behind.
$ws.Rows.Item(147).Insert()

# Populate the newly-inserted row 147 with the record.
$ws.Cells.Item(147, 1).Value = 10
$ws.Cells.Item(147, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(147, 3).Value = "La Araucanía"
$ws.Cells.Item(147, 4).Value = 44978
$ws.Cells.Item(147, 5).Value = 9
$ws.Cells.Item(147, 6).Value = 100114007
$ws.Cells.Item(147, 7).Value = "Jengibre"
$ws.Cells.Item(147, 8).Value = "Sin especificar"
$ws.Cells.Item(147, 9).Value = "Primera"
$ws.Cells.Item(147, 10).Value = 15
$ws.Cells.Item(147, 11).Value = 35000
$ws.Cells.Item(147, 12).Value = 35000
$ws.Cells.Item(147, 13).Value = 35000
$ws.Cells.Item(147, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(147, 15).Value = "Perú"
$ws.Cells.Item(147, 16).Value = 2692
$ws.Cells.Item(147, 17).Value = 13
$ws.Cells.Item(147, 18).Value = "Hortaliza"

# Match the date-number-format style used by the other rows' "Fecha" column.
$ws.Cells.Item(147, 4).NumberFormat = $ws.Cells.Item(148, 4).NumberFormat
